$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "data as of" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Julio de 2020 a las 17:51"

# Row 4
$ws.Range("B4").Value = 4448504
$ws.Range("C4").Value = 15094
$ws.Range("D4").Value = 2139014
$ws.Range("E4").Value = 2158038
$ws.Range("G4").Value = 377
$ws.Range("H4").Value = 151452

# Row 6
$ws.Range("B6").Value = 1516738
$ws.Range("C6").Value = 34235
$ws.Range("D6").Value = 971330
$ws.Range("E6").Value = 511542
$ws.Range("G6").Value = 418
$ws.Range("H6").Value = 33866

# Row 13
$ws.Range("B13").Value = 300692
$ws.Range("C13").Value = 581
$ws.Range("G13").Value = 119
$ws.Range("H13").Value = 45878

# Row 18
$ws.Range("D18").Value = 198756
$ws.Range("E18").Value = 12418

# Row 21
$ws.Range("B21").Value = 207487
$ws.Range("C21").Value = 108
$ws.Range("E21").Value = 6881

# Row 24
$ws.Range("A24").Value = "Irak"
$ws.Range("B24").Value = 115332
$ws.Range("C24").Value = 2747
$ws.Range("D24").Value = 81062
$ws.Range("E24").Value = 29735
$ws.Range("G24").Value = 77
$ws.Range("H24").Value = 4535

# Row 25
$ws.Range("A25").Value = "Canada"
$ws.Range("B25").Value = 114597
$ws.Range("D25").Value = 99860
$ws.Range("E25").Value = 5836
$ws.Range("H25").Value = 8901

# Row 26
$ws.Range("B26").Value = 109880
$ws.Range("C26").Value = 283
$ws.Range("D26").Value = 106603
$ws.Range("E26").Value = 3110
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 167

# Row 45
$ws.Range("D45").Value = 45893
$ws.Range("E45").Value = 5277

# Row 72
$ws.Range("A72").Value = "Chequia"
$ws.Range("B72").Value = 15684
$ws.Range("C72").Value = 168
$ws.Range("D72").Value = 11428
$ws.Range("E72").Value = 3882
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 374

# Row 73
$ws.Range("A73").Value = "Costa de Marfil"
$ws.Range("B73").Value = 15655
$ws.Range("D73").Value = 10361
$ws.Range("E73").Value = 5198
$ws.Range("H73").Value = 96

# Row 76
$ws.Range("B76").Value = 15200
$ws.Range("C76").Value = 653
$ws.Range("D76").Value = 6526
$ws.Range("E76").Value = 8435
$ws.Range("G76").Value = 11
$ws.Range("H76").Value = 239

# Row 95
$ws.Range("B95").Value = 6375
$ws.Range("C95").Value = 54
$ws.Range("D95").Value = 4855
$ws.Range("E95").Value = 1407
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 113

# Row 99
$ws.Range("A99").Value = "Albania"
$ws.Range("B99").Value = 4997
$ws.Range("C99").Value = 117
$ws.Range("D99").Value = 2789
$ws.Range("E99").Value = 2060
$ws.Range("G99").Value = 4
$ws.Range("H99").Value = 148

# Row 100
$ws.Range("A100").Value = "Croacia"
$ws.Range("B100").Value = 4923
$ws.Range("C100").Value = 42
$ws.Range("D100").Value = 4034
$ws.Range("E100").Value = 749
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 140

# Row 141
$ws.Range("B141").Value = 1182
$ws.Range("C141").Value = 6
$ws.Range("D141").Value = 1042
$ws.Range("E141").Value = 129

# Row 192
$ws.Range("B192").Value = 63
$ws.Range("C192").Value = 1
$ws.Range("E192").Value = 52
